$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Developer name
$ws.Range("C3").Value = "Nishant Malhotra"

# Test case 1 (row 7)
$ws.Range("F7").Value = "789, 1010, 1000.0, 2020-01-01, 2.50"
$ws.Range("G7").Value = "account_number is 789, management_fees is 2.50"

# Test case 2 (row 8)
$ws.Range("F8").Value = 'management_fees = "invalid"'
$ws.Range("G8").Value = "management_fee defaults to 2.55"

# Test case 3 (row 9)
$ws.Range("F9").Value = "date_created = 2020-01-01"
$ws.Range("G9").Value = "Returns 0.50 (Management fee waived)"

# Test case 4 (row 10)
$ws.Range("F10").Value = "date_created = Today - 10 years"
$ws.Range("G10").Value = "Returns 0.50 (Management fee waived)"

# Test case 5 (row 11)
$ws.Range("F11").Value = "date_created = Today, fee = 2.50"
$ws.Range("G11").Value = "Returns 3.00 (0.50 base + 2.50 fee)"

# Test case 6 (row 12)
$ws.Range("F12").Value = "date_created = 2010-01-01"
$ws.Range("G12").Value = 'String contains "Management Fees: Waived"'

# Test case 7 (row 13)
$ws.Range("F13").Value = "date_created = Today, fee = 2.50"
$ws.Range("G13").Value = 'String contains "Management Fees: 2.50"'

# Update the view to match where the author last left the cursor/selection
$ws.Range("G13").Select()
